# Update the cryptos list worksheet with freshly scraped values.
# Cells whose new text looks like a plain number (e.g. "1.0000", "241.92")
# are written via Formula with a leading apostrophe so Excel stores them
# as literal text (matching the original inline-string cell type) instead
# of silently converting them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '30.008.92'
$ws.Range("E2").Value = '  -1.06%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.900.73'
$ws.Range("E3").Value = '  -1.59%  '

# Row 4 - TetherUSD
$ws.Range("D4").Formula = "'1.0000"
$ws.Range("E4").Value = '  -0.09%  '

# Row 5 - XRP
$ws.Range("D5").Formula = "'0.7409"
$ws.Range("E5").Value = '  -0.86%  '

# Row 6 - BNB
$ws.Range("D6").Formula = "'241.92"
$ws.Range("E6").Value = '  -0.46%  '

# Row 7 - USDC
$ws.Range("D7").Formula = "'1.000"

# Row 8 - Cardano
$ws.Range("D8").Formula = "'0.3060"
$ws.Range("E8").Value = '  -3.55%  '

# Row 9 - Solana
$ws.Range("D9").Formula = "'25.82"
$ws.Range("E9").Value = '  -6.34%  '

# Row 10 - Dogecoin
$ws.Range("D10").Formula = "'0.06888"
$ws.Range("E10").Value = '  -3.32%  '

# Row 11 - TRON
$ws.Range("D11").Formula = "'0.08010"
$ws.Range("E11").Value = '  -0.53%  '

# Row 12 - Polygon
$ws.Range("D12").Formula = "'0.7557"
$ws.Range("E12").Value = '  -2.80%  '

# Row 13 - WrappedEther
$ws.Range("D13").Value = '1.918.93'
$ws.Range("E13").Value = '  -0.97%  '

# Row 14 - Polkadot
$ws.Range("D14").Formula = "'5.228"
$ws.Range("E14").Value = '  -3.19%  '

# Row 15 - Litecoin
$ws.Range("D15").Formula = "'91.10"
$ws.Range("E15").Value = '  -2.22%  '

# Row 16 - Uniswap
$ws.Range("D16").Formula = "'6.132"
$ws.Range("E16").Value = '  +1.93%  '

# Row 17 - WrappedBTC
$ws.Range("D17").Value = '29.988.79'
$ws.Range("E17").Value = '  -1.08%  '

# Row 18 - Avalanche
$ws.Range("E18").Value = '  -4.14%  '

# Row 19 - ShibaInu
$ws.Range("D19").Formula = "'0.000007741"
$ws.Range("E19").Value = '  -2.13%  '

# Row 20 - BitcoinCash
$ws.Range("D20").Formula = "'235.89"
$ws.Range("E20").Value = '  -6.30%  '

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = '2.151.84'
$ws.Range("E22").Value = '  -0.40%  '

# Row 23 - BinanceUSD
$ws.Range("D23").Formula = "'0.9999"
$ws.Range("E23").Value = '  -0.05%  '

# Row 24 - Chainlink
$ws.Range("D24").Formula = "'7.044"
$ws.Range("E24").Value = '  +5.68%  '

# Row 25 - Monero
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Formula = "'9.289"
$ws.Range("E25").Value = '  -2.91%  '

# Row 26 - Cosmos
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Formula = "'166.92"
$ws.Range("E26").Value = '  +0.89%  '

# Row 27 - EthereumClassic
$ws.Range("D27").Formula = "'18.78"
$ws.Range("E27").Value = '  -1.57%  '

# Row 28 - Stellar
$ws.Range("D28").Formula = "'0.1262"
$ws.Range("E28").Value = '  -2.61%  '

# Row 29 - LidoDAOToken
$ws.Range("D29").Formula = "'2.033"
$ws.Range("E29").Value = '  -6.73%  '

# Row 30 - Toncoin
$ws.Range("D30").Formula = "'1.360"
$ws.Range("E30").Value = '  -0.56%  '

# Row 31 - PancakeSwap
$ws.Range("D31").Formula = "'1.529"
$ws.Range("E31").Value = '  -1.81%  '

# Row 32 - Filecoin
$ws.Range("D32").Formula = "'4.289"
$ws.Range("E32").Value = '  -3.00%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Formula = "'4.029"
$ws.Range("E33").Value = '  -2.76%  '

# Row 34 - Hedera
$ws.Range("D34").Formula = "'0.05283"
$ws.Range("E34").Value = '  +0.73%  '

# Row 35 - ARBITRUM
$ws.Range("D35").Formula = "'1.278"
$ws.Range("E35").Value = '  -2.93%  '

# Row 36 - ImmutableX
$ws.Range("D36").Formula = "'0.7379"
$ws.Range("E36").Value = '  -2.63%  '

# Row 37 - HuobiToken
$ws.Range("D37").Formula = "'2.720"
$ws.Range("E37").Value = '  -2.37%  '

# Row 38 - VeChain
$ws.Range("D38").Formula = "'0.01933"
$ws.Range("E38").Value = '  -1.07%  '

# Row 39 - MXToken
$ws.Range("D39").Formula = "'2.763"
$ws.Range("E39").Value = '  -1.28%  '

# Row 40 - FraxShare
$ws.Range("D40").Formula = "'6.226"
$ws.Range("E40").Value = '  -4.54%  '

# Row 41 - TheSandbox
$ws.Range("D41").Formula = "'0.4440"
$ws.Range("E41").Value = '  -2.01%  '

# Row 42 - Aave
$ws.Range("E42").Value = '  -6.61%  '

# Row 43 - RenderToken
$ws.Range("D43").Formula = "'1.946"
$ws.Range("E43").Value = '  -1.32%  '

# Row 44 - PaxDollar
$ws.Range("E44").Value = '  -0.02%  '

# Row 45 - TrustWalletToken
$ws.Range("D45").Formula = "'0.8319"
$ws.Range("E45").Value = '  -1.29%  '

# Row 46 - Aptos
$ws.Range("D46").Formula = "'7.653"
$ws.Range("E46").Value = '  -0.79%  '

# Row 47 - Quant
$ws.Range("D47").Formula = "'100.95"
$ws.Range("E47").Value = '  -0.79%  '

# Row 48 - EnergySwap
$ws.Range("D48").Formula = "'9.757"
$ws.Range("E48").Value = '  -2.25%  '

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = '2.056.31'
$ws.Range("E49").Value = '  -0.71%  '

# Row 50 - Elrond
$ws.Range("D50").Formula = "'36.53"
$ws.Range("E50").Value = '  -3.20%  '

# Row 51 - Algorand
$ws.Range("D51").Formula = "'0.1164"
$ws.Range("E51").Value = '  -4.30%  '
